$d = $word.ActiveDocument
$nbsp = [char]0x00A0

# 1. Title: "Parking Garage Case Study" -> "The Lake Problem"
$d.Content.Find.Execute("Parking", $false, $false, $false, $false, $false, $true, 1, $false, "The", 2)
$d.Content.Find.Execute("Garage", $false, $false, $false, $false, $false, $true, 1, $false, "Lake", 2)
$d.Content.Find.Execute("Case Study", $false, $false, $false, $false, $false, $true, 1, $false, "Problem", 2)

# 2. Remove the Subtitle paragraph ("Real Options") entirely.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Style.NameLocal -eq "Subtitle") {
        $p.Range.Delete()
        break
    }
}

# 3. Date paragraph: "Fri., Mar. 8" -> "Fri., Mar. 22" (keep the non-breaking space before the day).
$d.Content.Find.Execute("Mar. 8", $false, $false, $false, $false, $false, $true, 1, $false, "Mar.${nbsp}22", 2)

Write-Output "Final paragraphs:"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Output "Para $i : style=$($p.Style.NameLocal) text=$($p.Range.Text)"
}
